# Add a new instruction row to the protected "Instructions" sheet.
#
# The "Instructions" sheet is protected (sheetProtection), so directly
# assigning a .Value to one of its cells is silently blocked. Copy/Paste
# operations, however, are still permitted while the sheet is protected.
# So we stage the new text on an unprotected sheet, insert the new (blank)
# row on the Instructions sheet, copy the staged cell into place, and then
# clean up the scratch cell.

$wb = $excel.ActiveWorkbook

$instructions = $wb.Worksheets.Item("Instructions")
$scratch = $wb.Worksheets.Item("ZEVs Supplied")

$newText = "(4) Dates must be on or after January 2nd, 2018."

# Stage the new string value far away from any real data on a sheet that
# isn't protected.
$stagingCell = $scratch.Cells.Item(100, 1)
$stagingCell.Value = $newText

# Insert a new row 13 (pushing the old row 13 - "(4) No more than 2000
# records may be entered." - down to row 14) on the Instructions sheet.
$instructions.Rows.Item(13).Insert()

# Copy the staged text into the newly inserted row.
$stagingCell.Copy($instructions.Cells.Item(13, 1))

# Remove the scratch value so the helper sheet is left unchanged.
$stagingCell.Clear()
